$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.413144
$ws.Range("H2").Value = 1.239432
$ws.Range("I2").Value = 0.4553782032534783
$ws.Range("J2").Value = 0.4553782032534782
$ws.Range("M2").Value = 0.8317113333333332
$ws.Range("N2").Value = 2.495134
$ws.Range("O2").Value = 0.0263454906755698
$ws.Range("P2").Value = 0.0263454906755698
$ws.Range("Q2").Value = 0.3436165470986666
$ws.Range("R2").Value = 3.092548923888
$ws.Range("S2").Value = 0.01199716220767224
$ws.Range("T2").Value = 0.01199716220767224

# Row 3
$ws.Range("G3").Value = 0.413144
$ws.Range("H3").Value = 1.239432
$ws.Range("I3").Value = 0.4553782032534783
$ws.Range("J3").Value = 0.4553782032534782
$ws.Range("O3").Value = 0.6529848313028861
$ws.Range("P3").Value = 0.6529848313028862
$ws.Range("Q3").Value = 8.516690609530666
$ws.Range("R3").Value = 76.650215485776
$ws.Range("S3").Value = 0.2973550592304839
$ws.Range("T3").Value = 0.2973550592304839

# Row 4
$ws.Range("G4").Value = 0.413144
$ws.Range("H4").Value = 1.239432
$ws.Range("I4").Value = 0.4553782032534783
$ws.Range("J4").Value = 0.4553782032534782
$ws.Range("M4").Value = 10.12334933333333
$ws.Range("N4").Value = 30.370048
$ws.Range("O4").Value = 0.3206696780215441
$ws.Range("P4").Value = 0.3206696780215441
$ws.Range("Q4").Value = 4.182401036970667
$ws.Range("R4").Value = 37.641609332736
$ws.Range("S4").Value = 0.1460259818153221
$ws.Range("T4").Value = 0.1460259818153221

# Row 5
$ws.Range("I5").Value = 0.3895918235379703
$ws.Range("J5").Value = 0.3895918235379702
$ws.Range("M5").Value = 0.8317113333333332
$ws.Range("N5").Value = 2.495134
$ws.Range("O5").Value = 0.0263454906755698
$ws.Range("P5").Value = 0.0263454906755698
$ws.Range("Q5").Value = 0.2939758561686667
$ws.Range("R5").Value = 2.645782705518
$ws.Range("S5").Value = 0.01026398775429783
$ws.Range("T5").Value = 0.01026398775429783

# Row 6
$ws.Range("I6").Value = 0.3895918235379703
$ws.Range("J6").Value = 0.3895918235379702
$ws.Range("O6").Value = 0.6529848313028861
$ws.Range("P6").Value = 0.6529848313028862
$ws.Range("Q6").Value = 7.286323766420667
$ws.Range("R6").Value = 65.57691389778601
$ws.Range("S6").Value = 0.2543975511699253
$ws.Range("T6").Value = 0.2543975511699253

# Row 7
$ws.Range("I7").Value = 0.3895918235379703
$ws.Range("J7").Value = 0.3895918235379702
$ws.Range("M7").Value = 10.12334933333333
$ws.Range("N7").Value = 30.370048
$ws.Range("O7").Value = 0.3206696780215441
$ws.Range("P7").Value = 0.3206696780215441
$ws.Range("Q7").Value = 3.578188932010667
$ws.Range("R7").Value = 32.203700388096
$ws.Range("S7").Value = 0.1249302846137471
$ws.Range("T7").Value = 0.1249302846137471

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1406516666666666
$ws.Range("H8").Value = 0.421955
$ws.Range("I8").Value = 0.1550299732085515
$ws.Range("J8").Value = 0.1550299732085515
$ws.Range("M8").Value = 0.8317113333333332
$ws.Range("N8").Value = 2.495134
$ws.Range("O8").Value = 0.0263454906755698
$ws.Range("P8").Value = 0.0263454906755698
$ws.Range("Q8").Value = 0.1169815852188889
$ws.Range("R8").Value = 1.05283426697
$ws.Range("S8").Value = 0.004084340713599729
$ws.Range("T8").Value = 0.00408434071359973

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1406516666666666
$ws.Range("H9").Value = 0.421955
$ws.Range("I9").Value = 0.1550299732085515
$ws.Range("J9").Value = 0.1550299732085515
$ws.Range("O9").Value = 0.6529848313028861
$ws.Range("P9").Value = 0.6529848313028862
$ws.Range("Q9").Value = 2.899441184465555
$ws.Range("R9").Value = 26.09497066019
$ws.Range("S9").Value = 0.101232220902477
$ws.Range("T9").Value = 0.101232220902477

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1406516666666666
$ws.Range("H10").Value = 0.421955
$ws.Range("I10").Value = 0.1550299732085515
$ws.Range("J10").Value = 0.1550299732085515
$ws.Range("M10").Value = 10.12334933333333
$ws.Range("N10").Value = 30.370048
$ws.Range("O10").Value = 0.3206696780215441
$ws.Range("P10").Value = 0.3206696780215441
$ws.Range("Q10").Value = 1.423865955982222
$ws.Range("R10").Value = 12.81479360384
$ws.Range("S10").Value = 0.04971341159247481
$ws.Range("T10").Value = 0.04971341159247482
